# Apply the "Detecting existing inst. names" edit:
#  - Column A/B swap meaning (A becomes c_inst_code, B becomes c_inst_name_hz)
#  - New column headers: A=c_inst_code, B=c_inst_name_hz, C=c_inst_name_code,
#    D=c_inst_type_code, E=c_inst_begin_dy, F=c_inst_addr_id, G=c_source
#  - Existing row of data (inst_code 3920 / 建初寺) keeps its position but gets
#    updated / reordered values
#  - Seven new data rows are appended (rows 3-9) for additional institutions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) stay the same text, just confirm/re-assert them in the new
# column order actually already matches existing shared strings positions.
$ws.Range("A1").Value = "c_inst_code"
$ws.Range("B1").Value = "c_inst_name_hz"
$ws.Range("C1").Value = "c_inst_name_code"
$ws.Range("D1").Value = "c_inst_type_code"
$ws.Range("E1").Value = "c_inst_begin_dy"
$ws.Range("F1").Value = "c_inst_addr_id"
$ws.Range("G1").Value = "c_source"

# Data rows 2-9
$data = @(
    @(3920, "建初寺", 777, "2", "19", "4540", "18417"),
    @(3921, "國慶寺", 2539, "2", "19", "5402", "18417"),
    @(3922, "南巖寺", 2540, "2", "19", "5404", "18417"),
    @(3923, "天童寺", 328, "2", "19", "5405", "18417"),
    @(3924, "天童寺", 328, "2", "20", "7627", "18417"),
    @(3925, "平陽寺", 2543, "2", "20", "5398", "18417"),
    @(3926, "天寧寺", 326, "2", "20", "7569", "18417"),
    @(3927, "虞山書院", 2350, "1", "20", "7546", "65006")
)

# Columns D:G hold text-typed numeric codes (e.g. "2", "19", "4540", "18417")
# that must be stored as shared-string text, not as numbers. Temporarily mark
# the range as Text before writing so Excel doesn't auto-coerce the values,
# then restore the default "Normal" cell style once all values are in place.
$textRange = $ws.Range("D2:G9")
$textRange.NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $row++
}

$textRange.Style = "Normal"
